$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (column C) date from 45192 to 45202 for all existing data rows (2-455)
$ws.Range("C2:C455").Value = 45202

# 2. Row 455 gains an explicit row height (matches default 15, but becomes "custom")
$ws.Rows.Item(455).RowHeight = 15

# 3. Add new row 456
$ws.Cells.Item(456, 1).Value = "A 46224-2023"
$ws.Cells.Item(456, 2).Value = 45196
$ws.Cells.Item(456, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(456, 3).Value = 45202
$ws.Cells.Item(456, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(456, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(456, 5).Value = "ÅSELE"
$ws.Cells.Item(456, 6).Value = "Naturvårdsverket"
$ws.Cells.Item(456, 7).Value = 3
$ws.Cells.Item(456, 8).Value = 0
$ws.Cells.Item(456, 9).Value = 0
$ws.Cells.Item(456, 10).Value = 0
$ws.Cells.Item(456, 11).Value = 0
$ws.Cells.Item(456, 12).Value = 0
$ws.Cells.Item(456, 13).Value = 0
$ws.Cells.Item(456, 14).Value = 0
$ws.Cells.Item(456, 15).Value = 0
$ws.Cells.Item(456, 16).Value = 0
$ws.Cells.Item(456, 17).Value = 0
$ws.Cells.Item(456, 18).WrapText = $true
$ws.Rows.Item(456).RowHeight = 15

# 4. Add new row 457
$ws.Cells.Item(457, 1).Value = "A 46693-2023"
$ws.Cells.Item(457, 2).Value = 45198
$ws.Cells.Item(457, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(457, 3).Value = 45202
$ws.Cells.Item(457, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(457, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(457, 5).Value = "ÅSELE"
$ws.Cells.Item(457, 6).Value = "Sveaskog"
$ws.Cells.Item(457, 7).Value = 14.8
$ws.Cells.Item(457, 8).Value = 0
$ws.Cells.Item(457, 9).Value = 0
$ws.Cells.Item(457, 10).Value = 0
$ws.Cells.Item(457, 11).Value = 0
$ws.Cells.Item(457, 12).Value = 0
$ws.Cells.Item(457, 13).Value = 0
$ws.Cells.Item(457, 14).Value = 0
$ws.Cells.Item(457, 15).Value = 0
$ws.Cells.Item(457, 16).Value = 0
$ws.Cells.Item(457, 17).Value = 0
$ws.Cells.Item(457, 18).WrapText = $true
